# Insert a new weekly price-record row above the current row 630 in the
# "Mango" subset sheet. Excel's native Rows.Insert() shifts the existing
# rows 630..706 down to 631..707 (carrying their values/formatting with
# them), which is exactly the row-shift visible across the whole diff.
# We then fill the freshly inserted row 630 with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(630).Insert()

$row = 630

$ws.Cells.Item($row, 1).Value  = 9
$ws.Cells.Item($row, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item($row, 3).Value  = "Metropolitana"
$ws.Cells.Item($row, 4).Value  = 45124
$ws.Cells.Item($row, 5).Value  = 13
$ws.Cells.Item($row, 6).Value  = "Fruta"
$ws.Cells.Item($row, 7).Value  = 100108
$ws.Cells.Item($row, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item($row, 9).Value  = 100108002
$ws.Cells.Item($row, 10).Value = "Mango"
$ws.Cells.Item($row, 11).Value = "Sin especificar"
$ws.Cells.Item($row, 12).Value = "Primera"
$ws.Cells.Item($row, 13).Value = 440
$ws.Cells.Item($row, 14).Value = 8500
$ws.Cells.Item($row, 15).Value = 8500
$ws.Cells.Item($row, 16).Value = 8500
$ws.Cells.Item($row, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item($row, 18).Value = "Brasil"
$ws.Cells.Item($row, 19).Value = 2125
$ws.Cells.Item($row, 20).Value = 4
